$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 and 5 (they are removed entirely in the new version)
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()

# Update row 2
$ws.Range("D2").Value = "MuSCs"
$ws.Range("G2").Value = 0.06082199999999999
$ws.Range("I2").Value = 0.1716860072883705
$ws.Range("J2").Value = 0.1716860072883705
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.6684413333333333
$ws.Range("N2").Value = 2.005324
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.040655938776
$ws.Range("R2").Value = 0.365903448984
$ws.Range("S2").Value = 0.1716860072883705
$ws.Range("T2").Value = 0.1716860072883705

# Update row 3
$ws.Range("A3").Value = "MuSCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.293441
$ws.Range("H3").Value = 0.880323
$ws.Range("I3").Value = 0.8283139927116295
$ws.Range("J3").Value = 0.8283139927116295
$ws.Range("M3").Value = 0.6684413333333333
$ws.Range("N3").Value = 2.005324
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.1961480932946667
$ws.Range("R3").Value = 1.765332839652
$ws.Range("S3").Value = 0.8283139927116295
$ws.Range("T3").Value = 0.8283139927116295
